$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the existing row 243, shifting rows 243:251
# down to 244:252 (their contents are unchanged by this operation).
$ws.Rows.Item(243).Insert()

# Populate the newly inserted row 243 with the new weekly price record.
$ws.Cells.Item(243, 1).Value = 11
$ws.Cells.Item(243, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(243, 3).Value = "Bíobío"
$ws.Cells.Item(243, 4).Value = 44714
$ws.Cells.Item(243, 5).Value = 8
$ws.Cells.Item(243, 6).Value = 100114001
$ws.Cells.Item(243, 7).Value = "Papa"
$ws.Cells.Item(243, 8).Value = "Asterix"
$ws.Cells.Item(243, 9).Value = "1a (cosecha)"
$ws.Cells.Item(243, 10).Value = 220
$ws.Cells.Item(243, 11).Value = 7000
$ws.Cells.Item(243, 12).Value = 8000
$ws.Cells.Item(243, 13).Value = 7545
$ws.Cells.Item(243, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(243, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(243, 16).Value = 302
$ws.Cells.Item(243, 17).Value = 25
$ws.Cells.Item(243, 18).Value = "Hortaliza"
